# Order upgrade - second try at reformatting - add options for order
#
# The "Bought 300 BNPQY @ 31.18" trade row (row 21) is being removed from
# the trade log entirely; every row below it shifts up by one to close the
# gap (the very last helper row that used to be row 63 disappears too).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole row 21 (BNPQY purchase) - shifts everything below up by one.
$ws.Rows("21:21").Delete()

# Leave the cursor where the author left it after the edit.
$ws.Range("F31").Select()
